$d = $word.ActiveDocument

# --- Part 1: merge the run split caused by the old "_GoBack" bookmark ---
# Remove the existing _GoBack bookmark (it sat in the middle of a sentence,
# splitting one logical run into two <w:r> elements).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Re-run Find/Replace over the now-unbroken sentence so Word collapses the
# (still split) runs back into a single run with uniform formatting.
$sentence = " csapatépítő oldalt készítettünk, amely lehetőséget biztosít a já" + `
    "tékosok számára, hogy mélyebben megértsék a játék mechanikáját és a csa" + `
    "patépítés stratégiáit. Az oldal bemutatja a különböző hősöket, azok ké" + `
    "pességeit, valamint a portálokat, amelyeket a játék elején, százalékos" + `
    " alapon sorsol ki. Emellett bemutatjuk a kiegészítőket, amelyeket a já" + `
    "tékos a mérkőzés során választhat, és részletes információkat adunk a" + `
    " játékban elérhető tárgyakról és azok kombinációiról. Mindezek segíte" + `
    "nek a játékosoknak abban, hogy hatékonyabb taktikákat alakítsanak ki " + `
    "és javítsák a teljesítményüket."

$find = $d.Content.Find
$find.Execute($sentence, $true, $false, $false, $false, $false, $true, 1, $false, $sentence, 2)

# --- Part 2: move the "Felhasználói dokumentáció" bullet ---
# It currently lives under the "Mózsik Emma:" heading; it belongs right
# after "Üveges Cintia:", before "Statikus oldalak fejlesztése".
$srcPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Felhasználói dokumentáció ") {
        $prev = $p.Previous()
        if ($prev -ne $null -and $prev.Range.Text.TrimEnd([char]13, [char]7) -eq "Mózsik Emma:") {
            $srcPara = $p
            break
        }
    }
}

if ($srcPara -eq $null) {
    throw "Could not locate the 'Felhasznaloi dokumentacio' paragraph under 'Mozsik Emma:'"
}

$srcRange = $srcPara.Range
$srcRange.Cut()

$destPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Statikus oldalak fejlesztése") {
        $destPara = $p
        break
    }
}

if ($destPara -eq $null) {
    throw "Could not locate the 'Statikus oldalak fejlesztese' paragraph"
}

$startPos = $destPara.Range.Start
$insertionPoint = $d.Range($startPos, $startPos)
$insertionPoint.Paste()

# --- Part 3: re-create the "_GoBack" bookmark at the new edit location ---
# (right after the moved run, mirroring where Word leaves it after a paste)
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Felhasználói dokumentáció ") {
        $prev = $p.Previous()
        if ($prev -ne $null -and $prev.Range.Text.TrimEnd([char]13, [char]7) -eq "Üveges Cintia: ") {
            $bmStart = $p.Range.End - 1
            $bmRange = $d.Range($bmStart, $bmStart)
            if ($d.Bookmarks.Exists("_GoBack")) {
                $d.Bookmarks.Item("_GoBack").Delete()
            }
            $d.Bookmarks.Add("_GoBack", $bmRange)
            break
        }
    }
}
